$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 140-189 (2025-05-30 batch) appended to the alert log.
# Format column A as text first so date-like strings ("2025-05-30") are
# stored as literal text (matching the rest of the sheet) instead of being
# auto-converted into a date serial by Excel's smart-entry parser.
$dateRange = $ws.Range("A140:A189")
$dateRange.NumberFormat = "@"

$ws.Range("A140").Value = '2025-05-30'
$ws.Range("B140").Value = 'Desconocido'
$ws.Range("C140").Value = 'Un aliado de la OTAN coge carrerilla con el avión ''invisible'' que rivaliza con el supercaza de EEUU'
$ws.Range("D140").Value = 'military'
$ws.Range("E140").Value = 2

$ws.Range("A141").Value = '2025-05-30'
$ws.Range("B141").Value = 'Desconocido'
$ws.Range("C141").Value = 'Un país clave de la OTAN pone en jaque a Rusia al revelar sus planes nucleares secretos'
$ws.Range("D141").Value = 'nuclear'
$ws.Range("E141").Value = 2

$ws.Range("A142").Value = '2025-05-30'
$ws.Range("B142").Value = 'Desconocido'
$ws.Range("C142").Value = 'Usan tierra de Fukushima en los canteros del primer ministro japonés para demostrar que no hay peligro'
$ws.Range("D142").Value = 'nuclear'
$ws.Range("E142").Value = 2

$ws.Range("A143").Value = '2025-05-30'
$ws.Range("B143").Value = 'Desconocido'
$ws.Range("C143").Value = 'Iberdrola asegura que en el apagón respondieron «siempre a los protocolos del Ministerio y de Redeia»'
$ws.Range("D143").Value = 'nuclear'
$ws.Range("E143").Value = 2

$ws.Range("A144").Value = '2025-05-30'
$ws.Range("B144").Value = 'Desconocido'
$ws.Range("C144").Value = 'Taco man'
$ws.Range("D144").Value = 'nuclear'
$ws.Range("E144").Value = 2

$ws.Range("A145").Value = '2025-05-30'
$ws.Range("B145").Value = 'Desconocido'
$ws.Range("C145").Value = 'Repsol llama a abandonar el "radicalismo ecologista" y reafirma su apuesta por los combustibles fósiles'
$ws.Range("D145").Value = 'nuclear'
$ws.Range("E145").Value = 2

$ws.Range("A146").Value = '2025-05-30'
$ws.Range("B146").Value = 'Desconocido'
$ws.Range("C146").Value = 'Galán señala a Aagesen y a Red Eléctrica en el apagón: "Iberdrola respondió según los protocolos"'
$ws.Range("D146").Value = 'nuclear'
$ws.Range("E146").Value = 2

$ws.Range("A147").Value = '2025-05-30'
$ws.Range("B147").Value = 'Desconocido'
$ws.Range("C147").Value = 'Zajárova: "Kiev quiere aumentar su valor con histeria pero solo le salen chichones"'
$ws.Range("D147").Value = 'nuclear'
$ws.Range("E147").Value = 2

$ws.Range("A148").Value = '2025-05-30'
$ws.Range("B148").Value = 'Desconocido'
$ws.Range("C148").Value = 'La fabricante de los misiles nucleares intercontinentales regresa a la Luna para competir contra SpaceX'
$ws.Range("D148").Value = 'nuclear'
$ws.Range("E148").Value = 2

$ws.Range("A149").Value = '2025-05-30'
$ws.Range("B149").Value = 'Desconocido'
$ws.Range("C149").Value = 'Atomfall se despide por todo lo alto con Wicked Isle, su última expansión que llega en junio con nuevos enemigos, armas y finales'
$ws.Range("D149").Value = 'nuclear'
$ws.Range("E149").Value = 2

$ws.Range("A150").Value = '2025-05-30'
$ws.Range("B150").Value = 'Desconocido'
$ws.Range("C150").Value = 'La guerra ya no es ficción: Algunos ya empezaron a desempolvar sus búnkeres y advierten que Europa vive una nueva era'
$ws.Range("D150").Value = 'nuclear'
$ws.Range("E150").Value = 2

$ws.Range("A151").Value = '2025-05-30'
$ws.Range("B151").Value = 'Desconocido'
$ws.Range("C151").Value = 'El CCIB se prepara para acoger los grandes congresos médicos mundiales del año'
$ws.Range("D151").Value = 'nuclear'
$ws.Range("E151").Value = 2

$ws.Range("A152").Value = '2025-05-30'
$ws.Range("B152").Value = 'Desconocido'
$ws.Range("C152").Value = 'Si algo no necesitaba la guerra de Ucrania era el tema “nuclear”. Rusia lo acaba de activar, literalmente'
$ws.Range("D152").Value = 'nuclear'
$ws.Range("E152").Value = 2

$ws.Range("A153").Value = '2025-05-30'
$ws.Range("B153").Value = 'Desconocido'
$ws.Range("C153").Value = 'La bandera fiscal a la que Mazón se aferra en medio de la tormenta'
$ws.Range("D153").Value = 'nuclear'
$ws.Range("E153").Value = 2

$ws.Range("A154").Value = '2025-05-30'
$ws.Range("B154").Value = 'Desconocido'
$ws.Range("C154").Value = 'Juegos para PC que corren en casi cualquier equipo'
$ws.Range("D154").Value = 'nuclear'
$ws.Range("E154").Value = 2

$ws.Range("A155").Value = '2025-05-30'
$ws.Range("B155").Value = 'Desconocido'
$ws.Range("C155").Value = 'Muere en un atentado en Rusia el militar al que Putin encargó el indiscriminado bombardeo de la ciudad de Mariúpol'
$ws.Range("D155").Value = 'nuclear'
$ws.Range("E155").Value = 2

$ws.Range("A156").Value = '2025-05-30'
$ws.Range("B156").Value = 'Desconocido'
$ws.Range("C156").Value = 'Por delante en tecnología, innovación y valores'
$ws.Range("D156").Value = 'nuclear'
$ws.Range("E156").Value = 2

$ws.Range("A157").Value = '2025-05-30'
$ws.Range("B157").Value = 'Desconocido'
$ws.Range("C157").Value = 'El vecino de España da luz verde al plan de modernización para albergar bombarderos nucleares'
$ws.Range("D157").Value = 'nuclear'
$ws.Range("E157").Value = 2

$ws.Range("A158").Value = '2025-05-30'
$ws.Range("B158").Value = 'Desconocido'
$ws.Range("C158").Value = '¿El inicio de una guerra con Trump?'
$ws.Range("D158").Value = 'nuclear'
$ws.Range("E158").Value = 2

$ws.Range("A159").Value = '2025-05-30'
$ws.Range("B159").Value = 'Desconocido'
$ws.Range("C159").Value = 'Momentos "Eureka" decisivos en la historia de la Física'
$ws.Range("D159").Value = 'nuclear'
$ws.Range("E159").Value = 2

$ws.Range("A160").Value = '2025-05-30'
$ws.Range("B160").Value = 'Desconocido'
$ws.Range("C160").Value = 'Una empresa europea planea generar 100 MW durante 40 años a partir de residuos nucleares'
$ws.Range("D160").Value = 'nuclear'
$ws.Range("E160").Value = 2

$ws.Range("A161").Value = '2025-05-30'
$ws.Range("B161").Value = 'Desconocido'
$ws.Range("C161").Value = 'Hell is Us: Lo hemos jugado y os contamos que esperar de está nueva IP'
$ws.Range("D161").Value = 'drone'
$ws.Range("E161").Value = 2

$ws.Range("A162").Value = '2025-05-30'
$ws.Range("B162").Value = 'Desconocido'
$ws.Range("C162").Value = 'Posibles anuncios de Summer Game Fest 2025: ¿Qué anunciarán compañías como Sony, Nintendo o Microsoft?'
$ws.Range("D162").Value = 'war'
$ws.Range("E162").Value = 2

$ws.Range("A163").Value = '2025-05-30'
$ws.Range("B163").Value = 'Desconocido'
$ws.Range("C163").Value = 'La temporada 4 de ''Bleach TYBW'' promete ser tremenda, pero tiene al propio Tite Kubo preocupado por el final del anime'
$ws.Range("D163").Value = 'war'
$ws.Range("E163").Value = 2

$ws.Range("A164").Value = '2025-05-30'
$ws.Range("B164").Value = 'Desconocido'
$ws.Range("C164").Value = '''Something beautiful'': el nuevo disco de Miley Cyrus suena a todo y a nada a la vez'
$ws.Range("D164").Value = 'war'
$ws.Range("E164").Value = 2

$ws.Range("A165").Value = '2025-05-30'
$ws.Range("B165").Value = 'Desconocido'
$ws.Range("C165").Value = 'La película de ‘Elden Ring’ quiere a esta estrella de Netflix y ‘Warfare’ como protagonista'
$ws.Range("D165").Value = 'war'
$ws.Range("E165").Value = 2

$ws.Range("A166").Value = '2025-05-30'
$ws.Range("B166").Value = 'Desconocido'
$ws.Range("C166").Value = 'Corcuera y Arabia: esta es la fascinante historia de las dos mujeres españolas que formaron parte del ejército americano durante la Segunda Guerra Mundial'
$ws.Range("D166").Value = 'war'
$ws.Range("E166").Value = 2

$ws.Range("A167").Value = '2025-05-30'
$ws.Range("B167").Value = 'Desconocido'
$ws.Range("C167").Value = '‘Something Beautiful’ de Miley Cyrus: algo bonito y seguro'
$ws.Range("D167").Value = 'war'
$ws.Range("E167").Value = 2

$ws.Range("A168").Value = '2025-05-30'
$ws.Range("B168").Value = 'Desconocido'
$ws.Range("C168").Value = 'Se filtra un nuevo juego que va camino de Nintendo Switch 2. Un espectáculo gráfico que nos tiene a muchos muy intrigados'
$ws.Range("D168").Value = 'war'
$ws.Range("E168").Value = 2

$ws.Range("A169").Value = '2025-05-30'
$ws.Range("B169").Value = 'Desconocido'
$ws.Range("C169").Value = 'El director de la película Elden Ring ya tiene algunas ideas para el reparto'
$ws.Range("D169").Value = 'war'
$ws.Range("E169").Value = 2

$ws.Range("A170").Value = '2025-05-30'
$ws.Range("B170").Value = 'Desconocido'
$ws.Range("C170").Value = 'Adiós al examen teórico de toda la vida: la DGT planea cambios radicales en el test con vídeos y respuestas mútiples'
$ws.Range("D170").Value = 'terror'
$ws.Range("E170").Value = 2

$ws.Range("A171").Value = '2025-05-30'
$ws.Range("B171").Value = 'Desconocido'
$ws.Range("C171").Value = '¡Estamos en problemas! Dos de las peores termitas del mundo se han cruzado en Florida.'
$ws.Range("D171").Value = 'terror'
$ws.Range("E171").Value = 2

$ws.Range("A172").Value = '2025-05-30'
$ws.Range("B172").Value = 'Desconocido'
$ws.Range("C172").Value = 'Habrá otra spinoff para TV de «The Batman»'
$ws.Range("D172").Value = 'terror'
$ws.Range("E172").Value = 2

$ws.Range("A173").Value = '2025-05-30'
$ws.Range("B173").Value = 'Desconocido'
$ws.Range("C173").Value = 'Una uruguaya que vive en España impresionada por el dinero que se paga por ir a una boda en nuestro país: «El verdadero terror debe ser...»'
$ws.Range("D173").Value = 'terror'
$ws.Range("E173").Value = 2

$ws.Range("A174").Value = '2025-05-30'
$ws.Range("B174").Value = 'Desconocido'
$ws.Range("C174").Value = '''Elden Ring'': Kit Connor suena para protagonizar el live-action del videojuego dirigido por Alex Garland'
$ws.Range("D174").Value = 'terror'
$ws.Range("E174").Value = 2

$ws.Range("A175").Value = '2025-05-30'
$ws.Range("B175").Value = 'Desconocido'
$ws.Range("C175").Value = 'Coca-Cola devuelve más agua de la que consume'
$ws.Range("D175").Value = 'terror'
$ws.Range("E175").Value = 2

$ws.Range("A176").Value = '2025-05-30'
$ws.Range("B176").Value = 'Desconocido'
$ws.Range("C176").Value = 'Otrora signos de terror, los cometas revelan ahora secretos del universo'
$ws.Range("D176").Value = 'terror'
$ws.Range("E176").Value = 2

$ws.Range("A177").Value = '2025-05-30'
$ws.Range("B177").Value = 'Desconocido'
$ws.Range("C177").Value = 'Top juegos indie imprescindibles para PC'
$ws.Range("D177").Value = 'terror'
$ws.Range("E177").Value = 2

$ws.Range("A178").Value = '2025-05-30'
$ws.Range("B178").Value = 'Desconocido'
$ws.Range("C178").Value = 'Five Nights at Freddy''s se va del streaming'
$ws.Range("D178").Value = 'terror'
$ws.Range("E178").Value = 2

$ws.Range("A179").Value = '2025-05-30'
$ws.Range("B179").Value = 'Desconocido'
$ws.Range("C179").Value = 'Planes familiares para el último fin de semana de mayo... y primero de junio'
$ws.Range("D179").Value = 'terror'
$ws.Range("E179").Value = 2

$ws.Range("A180").Value = '2025-05-30'
$ws.Range("B180").Value = 'Desconocido'
$ws.Range("C180").Value = '¿Está La calle del terror: La reina del baile basada en hechos reales?'
$ws.Range("D180").Value = 'terror'
$ws.Range("E180").Value = 2

$ws.Range("A181").Value = '2025-05-30'
$ws.Range("B181").Value = 'Desconocido'
$ws.Range("C181").Value = 'Primer piso de espanto'
$ws.Range("D181").Value = 'terror'
$ws.Range("E181").Value = 2

$ws.Range("A182").Value = '2025-05-30'
$ws.Range("B182").Value = 'Desconocido'
$ws.Range("C182").Value = 'No solo Grupo Fugitivo: El Cártel del Golfo arrastra un historial de matanzas desde hace más de 40 años'
$ws.Range("D182").Value = 'terror'
$ws.Range("E182").Value = 2

$ws.Range("A183").Value = '2025-05-30'
$ws.Range("B183").Value = 'Desconocido'
$ws.Range("C183").Value = 'Wes Anderson vuelve tan simétrico como siempre tras salir escaldado de Cannes. Lo mejor y lo peor de los estrenos de cine'
$ws.Range("D183").Value = 'terror'
$ws.Range("E183").Value = 2

$ws.Range("A184").Value = '2025-05-30'
$ws.Range("B184").Value = 'Desconocido'
$ws.Range("C184").Value = 'Todos los nuevos juegos que llegan a Xbox para cerrar mayo por todo lo alto'
$ws.Range("D184").Value = 'terror'
$ws.Range("E184").Value = 2

$ws.Range("A185").Value = '2025-05-30'
$ws.Range("B185").Value = 'Desconocido'
$ws.Range("C185").Value = 'Guía para disfrutar el Sundance Film Festival CDMX 2025: sedes, películas y boletos'
$ws.Range("D185").Value = 'terror'
$ws.Range("E185").Value = 2

$ws.Range("A186").Value = '2025-05-30'
$ws.Range("B186").Value = 'Desconocido'
$ws.Range("C186").Value = 'Por el honor y 330 mil dólares: Colo Colo cierra su Copa Libertadores más amarga venciendo al rústico Bucaramanga'
$ws.Range("D186").Value = 'terror'
$ws.Range("E186").Value = 2

$ws.Range("A187").Value = '2025-05-30'
$ws.Range("B187").Value = 'Desconocido'
$ws.Range("C187").Value = 'Estado Islámico ataca por primera vez a los sucesores de Bashar al Asad'
$ws.Range("D187").Value = 'terror'
$ws.Range("E187").Value = 2

$ws.Range("A188").Value = '2025-05-30'
$ws.Range("B188").Value = 'Desconocido'
$ws.Range("C188").Value = 'Ángel Di María, ante el desafío de romper una nueva pared: ahora, con Rosario Central'
$ws.Range("D188").Value = 'terror'
$ws.Range("E188").Value = 2

$ws.Range("A189").Value = '2025-05-30'
$ws.Range("B189").Value = 'Desconocido'
$ws.Range("C189").Value = 'Todos los juegos de lanzamiento de Switch 2 con los que podrás estrenar la consola de Nintendo'
$ws.Range("D189").Value = 'attack'
$ws.Range("E189").Value = 2

# Clear the temporary text-number-format so the cells end up unstyled,
# same as every other data row in the sheet.
$dateRange.Style = "Normal"

Write-Host "Added rows 140-189"
